$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 128
$ws.Range("I2").Value = 353
$ws.Range("J2").Value = 1484
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 393
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = 259
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 147
$ws.Range("T2").Value = 252
$ws.Range("U2").Value = 23
$ws.Range("V2").Value = 2231
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2286
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 37
$ws.Range("AA2").Value = 14
